$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header row (A1:D1) to short machine-friendly column names
$ws.Range("A1").Value2 = "mx_state"
$ws.Range("B1").Value2 = "mx_municipality"
$ws.Range("C1").Value2 = "n_matriculas"
$ws.Range("D1").Value2 = "pct_matriculas"

# 2. Title-case Spanish connector words ("de", "del", "de la", "de los", "el", "y")
#    in state/municipality names throughout the data, and normalize the
#    all-caps grand-total label "TOTAL" to "Total".
$updates = @(
    @("B7", "Pabellón De Arteaga"),
    @("B8", "Rincón De Romos"),
    @("B9", "San Francisco De Los Romo"),
    @("B21", "Amatenango Del Valle"),
    @("B23", "Bejucal De Ocampo"),
    @("B30", "Comitán De Domínguez"),
    @("B44", "Mazapa De Madero"),
    @("B49", "Ocozocoautla De Espinosa"),
    @("B56", "Salto De Agua"),
    @("B58", "San Cristóbal De Las Casas"),
    @("B86", "Guadalupe Y Calvo"),
    @("B88", "Hidalgo Del Parral"),
    @("B93", "San Francisco Del Oro"),
    @("A95", "Ciudad De México"),
    @("A109", "Coahuila De Zaragoza"),
    @("B114", "San Juan De Sabinas"),
    @("B123", "Villa De Álvarez"),
    @("B126", "Coneto De Comonfort"),
    @("B140", "San Juan Del Río"),
    @("A148", "Estado De México"),
    @("B148", "Acambay De Ruíz Castañeda"),
    @("B150", "Almoloya De Juárez"),
    @("B154", "Atizapán De Zaragoza"),
    @("B158", "Coacalco De Berriozábal"),
    @("B164", "Ecatepec De Morelos"),
    @("B170", "Naucalpan De Juárez"),
    @("B175", "San Felipe Del Progreso"),
    @("B176", "San José Del Rincón"),
    @("B185", "Tlalnepantla De Baz"),
    @("B188", "Valle De Bravo"),
    @("B190", "Villa Del Carbón"),
    @("B195", "Apaseo El Grande"),
    @("B201", "Dolores Hidalgo Cuna De La Independencia Nacional"),
    @("B205", "Jaral Del Progreso"),
    @("B210", "Purísima Del Rincón"),
    @("B216", "San Francisco Del Rincón"),
    @("B218", "San Luis De La Paz"),
    @("B219", "San Miguel De Allende"),
    @("B220", "Santa Cruz De Juventino Rosas"),
    @("B221", "Silao De La Victoria"),
    @("B226", "Valle De Santiago"),
    @("B232", "Acapulco De Juárez"),
    @("B233", "Alcozauca De Guerrero"),
    @("B237", "Atoyac De Álvarez"),
    @("B238", "Ayutla De Los Libres"),
    @("B240", "Buenavista De Cuéllar"),
    @("B241", "Chilpancingo De Los Bravo"),
    @("B242", "Coahuayutla De José María Izazaga"),
    @("B244", "Coyuca De Benítez"),
    @("B247", "Cuetzala Del Progreso"),
    @("B254", "Iguala De La Independencia"),
    @("B256", "Ixcateopan De Cuauhtémoc"),
    @("B262", "Mártir De Cuilapan"),
    @("B270", "Taxco De Alarcón"),
    @("B273", "Tepecoacuilco De Trujano"),
    @("B276", "Tlalixtaquilla De Maldonado"),
    @("B277", "Tlapa De Comonfort"),
    @("B279", "Técpan De Galeana"),
    @("B281", "Zihuatanejo De Azueta"),
    @("B285", "Agua Blanca De Iturbide"),
    @("B288", "Cuautepec De Hinojosa"),
    @("B289", "Huasca De Ocampo"),
    @("B290", "Huejutla De Reyes"),
    @("B296", "Mineral Del Monte"),
    @("B297", "Pachuca De Soto"),
    @("B300", "Progreso De Obregón"),
    @("B302", "Santiago Tulantepec De Lugo Guerrero"),
    @("B303", "Santiago De Anaya"),
    @("B305", "Tepehuacán De Guerrero"),
    @("B306", "Tula De Allende"),
    @("B307", "Tulancingo De Bravo"),
    @("B308", "Zacualtipán De Ángeles"),
    @("B315", "Atotonilco El Alto"),
    @("B316", "Autlán De Navarro"),
    @("B327", "Huejuquilla El Alto"),
    @("B328", "Ixtlahuacán De Los Membrillos"),
    @("B329", "Ixtlahuacán Del Río"),
    @("B333", "Lagos De Moreno"),
    @("B340", "Ojuelos De Jalisco"),
    @("B343", "San Cristóbal De La Barranca"),
    @("B344", "San Juan De Los Lagos"),
    @("B346", "San Miguel El Alto"),
    @("B347", "San Sebastián Del Oeste"),
    @("B349", "Tamazula De Gordiano"),
    @("B352", "Teocuitatlán De Corona"),
    @("B353", "Tepatitlán De Morelos"),
    @("B355", "Tizapán El Alto"),
    @("B358", "Unión De San Antonio"),
    @("B359", "Unión De Tula"),
    @("B361", "Yahualica De González Gallo"),
    @("B365", "Zapotlán El Grande"),
    @("A367", "Michoacán De Ocampo"),
    @("B376", "Coalcomán De Vázquez Pallares"),
    @("B378", "Cojumatlán De Régules"),
    @("B439", "Puente De Ixtla"),
    @("B442", "Tlaltizapán De Zapata"),
    @("B448", "Amatlán De Cañas"),
    @("B449", "Bahía De Banderas"),
    @("B467", "San Nicolás De Los Garza"),
    @("B470", "Acatlán De Pérez Figueroa"),
    @("B474", "Cuilápam De Guerrero"),
    @("B476", "Heroica Ciudad De Ejutla De Crespo"),
    @("B477", "Heroica Ciudad De Huajuapan De León"),
    @("B478", "Heroica Ciudad De Juchitán De Zaragoza"),
    @("B479", "Heroica Ciudad De Tlaxiaco"),
    @("B480", "Huajuapan De León"),
    @("B484", "Miahuatlán De Porfirio Díaz"),
    @("B485", "Oaxaca De Juárez"),
    @("B486", "Pinotepa De Don Luis"),
    @("B489", "San Antonino El Alto"),
    @("B491", "San Felipe Jalapa De Díaz"),
    @("B492", "San Francisco Del Mar"),
    @("B508", "San Miguel Del Puerto"),
    @("B510", "San Pablo Villa De Mitla"),
    @("B537", "Tamazulápam Del Espíritu Santo"),
    @("B538", "Teotitlán De Flores Magón"),
    @("B539", "Tezoatlán De Segura Y Luna"),
    @("B541", "Villa De Tututepec"),
    @("B542", "Zimatlán De Álvarez"),
    @("B557", "Cuetzalan Del Progreso"),
    @("B563", "Ixcamilpa De Guerrero"),
    @("B564", "Izúcar De Matamoros"),
    @("B568", "Los Reyes De Juárez"),
    @("B579", "San Salvador El Seco"),
    @("B586", "Tlacotepec De Benito Juárez"),
    @("B599", "Amealco De Bonfil"),
    @("B604", "Jalpan De Serra"),
    @("B608", "San Juan Del Río"),
    @("B625", "Santa María Del Río"),
    @("B629", "Tanquián De Escobedo"),
    @("B630", "Villa De Arriaga"),
    @("B631", "Villa De Guadalupe"),
    @("B632", "Villa De Ramos"),
    @("B653", "Nacozari De García"),
    @("B679", "Soto La Marina"),
    @("B682", "Acuamanala De Miguel Hidalgo"),
    @("B686", "Tetla De La Solidaridad"),
    @("A690", "Veracruz De Ignacio De La Llave"),
    @("B692", "Alto Lucero De Gutiérrez Barrios"),
    @("B693", "Amatlán De Los Reyes"),
    @("B697", "Boca Del Río"),
    @("B701", "Cazones De Herrera"),
    @("B708", "Cosamaloapan De Carpio"),
    @("B719", "Hueyapan De Ocampo"),
    @("B724", "Juchique De Ferrer"),
    @("B727", "Lerdo De Tejada"),
    @("B730", "Martínez De La Torre"),
    @("B735", "Mixtla De Altamirano"),
    @("B742", "Paso Del Macho"),
    @("B745", "Poza Rica De Hidalgo"),
    @("B751", "Sayula De Alemán"),
    @("B752", "Soledad De Doblado"),
    @("B765", "Vega De Alatorre"),
    @("B781", "Cañitas De Felipe Pescador"),
    @("B782", "Concepción Del Oro"),
    @("B783", "El Plateado De Joaquín Amaro"),
    @("B796", "Moyahua De Estrada"),
    @("B797", "Nochistlán De Mejía"),
    @("B798", "Noria De Ángeles"),
    @("B807", "Teúl De González Ortega"),
    @("B808", "Tlaltenango De Sánchez Román"),
    @("B812", "Villa De Cos"),
    @("A816", "Total")
)

foreach ($pair in $updates) {
    $ws.Range($pair[0]).Value2 = $pair[1]
}

# 3. Remove the trailing footnote/source rows (818-822) that are no longer
#    part of the clean dataset. Row 817 is already blank and acts as the
#    separator before this block.
$ws.Range("A818:D822").EntireRow.Delete()
